# Update validations in "add deals" test case data (deals sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("deals")

# --- New columns: O (predictedCloseDate) and P (actualCloseDate) ---
$ws.Range("O1").Value = "predictedCloseDate"
$ws.Range("P1").Value = "actualCloseDate"
# Match the yellow header fill used by the rest of row 1
$ws.Range("O1:P1").Interior.ColorIndex = 6

# Apply a date number format (this creates the numFmtId=15 style, and must be
# created before the text styles below so that it lands at cellXfs index 3)
$ws.Range("O2:P3").NumberFormat = "d-mmm-yy"
$ws.Range("O2").Value = 43432
$ws.Range("P2").Value = 43434
$ws.Range("O3").Value = 43429
$ws.Range("P3").Value = 43430

# Give the new columns a sensible display width
$ws.Columns.Item(15).ColumnWidth = 18.0221354166667
$ws.Columns.Item(16).ColumnWidth = 14.5924479166667

# --- probability column (E) switches from numeric to free-text values ---
# Setting number format to Text ("@") before assigning the values keeps them
# as text ("80"/"60") instead of being re-interpreted as numbers.
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E2:E3").NumberFormat = "@"
$ws.Range("E2").Value = "80"
$ws.Range("E3").Value = "60"

# --- Update the sheet's active selection to column F ---
$ws.Activate()
$ws.Range("F:F").Select()
